# Generate Report for Handoff
# Updates handoff timestamps and sets the handoff-priority marker ("ht")
# for the rows whose handback priority needed to match the handoff type.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$rows = @(7, 8, 9, 10, 12, 14)

foreach ($r in $rows) {
    # Overview sheet: column G = "Latest HO Xliff Generate Date"
    $wsOverview.Range("G$r").Value = "2016-08-27 06:19:53"

    # zh-cn sheet: column E = "Priority", column H = "Latest Handoff Datetime"
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-08-27 06:19:49"

    # de-de sheet: column E = "Priority", column H = "Latest Handoff Datetime"
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-08-27 06:19:53"
}
